# indexMap10x16 - "modified msh, and added pythonNoteBook"
#
# This script:
#  1) Grows the existing 10x16 index-map grid's row height (12.8 -> 27.35,
#     custom height) for the header/body rows (12-27).
#  2) Re-colors a few of the existing fills:
#       - the cyan/green header fill's background swatch (fg unchanged)
#       - the magenta "msh" block (rows 20-24, cols E:L) background swatch
#       - the magenta/yellow block (rows 20-24, cols M:T) recolored to a
#         brand-new green fill
#  3) Adds a small 4x4 legend / "python notebook" grid in rows 33-36
#     (cols C:F) using a new 14pt font, reusing + adding a couple of new
#     fills, with blank (style-only) corner cells and numeric values.
#  4) Updates the view (selection + scroll position) and row heights for
#     the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row heights for the existing grid (rows 12-27): 12.8 -> 27.35
# ---------------------------------------------------------------------
$ws.Range("A12:A27").RowHeight = 27.35

# ---------------------------------------------------------------------
# 2) Re-color existing fills (background swatch only unless noted)
# ---------------------------------------------------------------------

# Header rows (12-14) and footer rows (25-27), cols E:T -- light cyan fill,
# background swatch CCFFCC -> 99FF99 (foreground/visible color unchanged).
$hdr1 = $ws.Range("E12:T14")
$hdr1.Interior.PatternColor = 10092441   # bg 99FF99
$hdr1.Interior.Color        = 16777113   # fg 99FFFF (unchanged)

$hdr2 = $ws.Range("E25:T27")
$hdr2.Interior.PatternColor = 10092441   # bg 99FF99
$hdr2.Interior.Color        = 16777113   # fg 99FFFF (unchanged)

# Rows 20-24, cols E:L -- magenta fill, background swatch CC99FF -> FF9999
# (foreground/visible color unchanged).
$mid = $ws.Range("E20:L24")
$mid.Interior.PatternColor = 10066431    # bg FF9999
$mid.Interior.Color        = 16751103    # fg FF99FF (unchanged)

# Rows 20-24, cols M:T -- recolored entirely to a new green/cyan fill.
$midRight = $ws.Range("M20:T24")
$midRight.Interior.PatternColor = 16777113   # bg 99FFFF
$midRight.Interior.Color        = 10092441   # fg 99FF99

# ---------------------------------------------------------------------
# 3) New "python notebook" legend grid: rows 33-36, cols C:F
# ---------------------------------------------------------------------

# First stamp every cell in the block with the existing bordered/no-fill
# style (copy format from B12, which already uses the hair border with no
# fill), then bump the font size to 14pt and set per-cell fill + values.
$ws.Range("B12").Copy()
$legend = $ws.Range("C33:F36")
$legend.PasteSpecial(-4122)          # xlPasteFormats
$legend.Font.Size = 14

# Row 33: blank / -2 / -4 / blank
$ws.Range("D33").Value = -2
$ws.Range("E33").Value = -4

# Row 34: 3 / 2 / 3 / 2
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 2

# Row 35: 1 / 0 / 1 / 0
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = 0

# Row 36: blank / -3 / -1 / blank
$ws.Range("D36").Value = -3
$ws.Range("E36").Value = -1

# Fill colors for the legend grid.
# C/F columns (rows 34-35) -- light blue fill (same swatch as fill idx 3).
$blueCells = $ws.Range("C34:C35,F34:F35")
$blueCells.Interior.PatternColor = 16764108   # bg CCCCFF
$blueCells.Interior.Color        = 16115663   # fg CFE7F5

# D33,E33,D36,E36 -- green/cyan fill (same swatch as the re-colored header).
$greenCells = $ws.Range("D33,E33,D36,E36")
$greenCells.Interior.PatternColor = 10092441  # bg 99FF99
$greenCells.Interior.Color        = 16777113  # fg 99FFFF

# D34 -- new salmon/pink fill.
$ws.Range("D34").Interior.PatternColor = 16751103  # bg FF99FF
$ws.Range("D34").Interior.Color        = 10066431  # fg FF9999

# E34 -- yellow fill (same swatch as the existing yellow fill).
$ws.Range("E34").Interior.PatternColor = 65535     # bg FFFF00
$ws.Range("E34").Interior.Color        = 6750207   # fg FFFF66

# D35 -- new magenta/pink fill.
$ws.Range("D35").Interior.PatternColor = 10066431  # bg FF9999
$ws.Range("D35").Interior.Color        = 13395711  # fg FF66CC

# E35 -- new green fill (same swatch as the M20:T24 recolor).
$ws.Range("E35").Interior.PatternColor = 16777113  # bg 99FFFF
$ws.Range("E35").Interior.Color        = 10092441  # fg 99FF99

# ---------------------------------------------------------------------
# 4) Row heights for the new rows + blank spacer rows, view/selection
# ---------------------------------------------------------------------
$ws.Range("A32:A38").RowHeight = 27.35

$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("T39").Select()
